# Converts the "Chapter"-flavoured section head (plus H2O structural
# marker paragraphs: NodeStart/NodeEnd/HeadSeparator/HeadFieldSeparator/
# HeadEnd/ChapterSpacer) into the "Section"-flavoured head used by the
# annotations export, collapsing each removed run of marker paragraphs
# down to a single invisible-separator paragraph where the diff calls
# for one.
#
# Because deleting a paragraph's Range renumbers every later paragraph,
# the edits below are applied from the highest paragraph index to the
# lowest so that indices already queued for processing stay valid.

$d = $word.ActiveDocument

function Set-ParaStyle($index, $styleName) {
    $d.Paragraphs.Item($index).Style = $styleName
}

function Set-ParaStyleAndText($index, $styleName, $text) {
    $p = $d.Paragraphs.Item($index)
    $p.Style = $styleName
    $p.Range.Text = $text
}

function Remove-Para($index) {
    $d.Paragraphs.Item($index).Range.Delete()
}

# --- trailing block: Resource 1.3 footer + final NodeEnd -----------------
Set-ParaStyleAndText 52 "invisibleseparator" " "
Remove-Para 51
Remove-Para 49

# --- Resource 1.3 header (NodeEnd/NodeStart/HeadSeparator -> separator) --
Set-ParaStyleAndText 46 "invisibleseparator" " "
Remove-Para 45
Remove-Para 44

Remove-Para 42

# --- Resource 1.2 header (NodeEnd/NodeStart/HeadSeparator -> separator) --
Set-ParaStyleAndText 39 "invisibleseparator" " "
Remove-Para 38
Remove-Para 37

Remove-Para 24

# --- Resource 1.1 header (NodeEnd/NodeStart/HeadSeparator removed) -------
Remove-Para 21
Remove-Para 20
Remove-Para 19

# --- chapter/section head block ------------------------------------------
Set-ParaStyleAndText 10 "invisibleseparator" " "
Set-ParaStyle 9 "SectionHeadnote"
Remove-Para 8
Set-ParaStyle 7 "SectionSubtitle"
Remove-Para 6
Set-ParaStyle 5 "SectionTitle"
Set-ParaStyle 4 "SectionNumber"
Remove-Para 3
Remove-Para 2
Remove-Para 1

# --- remove now-unused custom paragraph styles from the style sheet ------
$stylesToRemove = @("ChapterSpacer", "HeadEnd", "HeadFieldSeparator", "HeadSeparator", "NodeEnd", "NodeStart")
foreach ($styleName in $stylesToRemove) {
    $d.Styles.Item($styleName).Delete()
}
